$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, $r, $a, $b, $c, $d, $e, $f) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
}

$ws = $wb.Worksheets.Item("PIR")
Set-RowValues $ws 18 '2026-02-01' '13:48:35' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 19 '2026-02-01' '13:48:36' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 20 '2026-02-01' '13:48:40' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 21 '2026-02-01' '13:48:41' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 22 '2026-02-01' '13:49:18' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 23 '2026-02-01' '13:49:18' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 24 '2026-02-01' '13:49:18' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 25 '2026-02-01' '13:49:19' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 26 '2026-02-01' '13:49:19' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 27 '2026-02-01' '13:49:19' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 28 '2026-02-01' '13:49:19' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 29 '2026-02-01' '13:49:19' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 30 '2026-02-01' '13:49:20' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 31 '2026-02-01' '13:49:20' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 32 '2026-02-01' '13:49:20' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 33 '2026-02-01' '13:49:20' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 34 '2026-02-01' '13:49:21' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 35 '2026-02-01' '13:49:21' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 36 '2026-02-01' '13:49:21' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 37 '2026-02-01' '13:49:22' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 38 '2026-02-01' '13:49:22' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 39 '2026-02-01' '13:49:22' '13:00' 'Bathroom' 'Motion Detected' 'Active'
Set-RowValues $ws 40 '2026-02-01' '13:49:23' '13:00' 'Bathroom' 'No Motion' 'Inactive'
Set-RowValues $ws 41 '2026-02-01' '13:49:24' '13:00' 'Bathroom' 'No Motion' 'Inactive'

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("E2").NumberFormat = "@"
Set-RowValues $ws 2 '2026-02-01' '13:48:33' '13:00' 'Bathroom' '90.0%' 'Active'
$ws.Range("E3").NumberFormat = "@"
Set-RowValues $ws 3 '2026-02-01' '13:48:37' '13:00' 'Bathroom' '89.1%' 'Active'
$ws.Range("E4").NumberFormat = "@"
Set-RowValues $ws 4 '2026-02-01' '13:48:38' '13:00' 'Bathroom' '86.2%' 'Active'
$ws.Range("E5").NumberFormat = "@"
Set-RowValues $ws 5 '2026-02-01' '13:49:18' '13:00' 'Bathroom' '84.4%' 'Active'
$ws.Range("E6").NumberFormat = "@"
Set-RowValues $ws 6 '2026-02-01' '13:49:18' '13:00' 'Bathroom' '82.7%' 'Active'
$ws.Range("E7").NumberFormat = "@"
Set-RowValues $ws 7 '2026-02-01' '13:49:18' '13:00' 'Bathroom' '83.6%' 'Active'
$ws.Range("E8").NumberFormat = "@"
Set-RowValues $ws 8 '2026-02-01' '13:49:19' '13:00' 'Bathroom' '97.1%' 'Active'
$ws.Range("E9").NumberFormat = "@"
Set-RowValues $ws 9 '2026-02-01' '13:49:19' '13:00' 'Bathroom' '96.6%' 'Active'
$ws.Range("E10").NumberFormat = "@"
Set-RowValues $ws 10 '2026-02-01' '13:49:20' '13:00' 'Bathroom' '90.3%' 'Active'
$ws.Range("E11").NumberFormat = "@"
Set-RowValues $ws 11 '2026-02-01' '13:49:20' '13:00' 'Bathroom' '87.4%' 'Active'
$ws.Range("E12").NumberFormat = "@"
Set-RowValues $ws 12 '2026-02-01' '13:49:21' '13:00' 'Bathroom' '86.4%' 'Active'
$ws.Range("E13").NumberFormat = "@"
Set-RowValues $ws 13 '2026-02-01' '13:49:21' '13:00' 'Bathroom' '84.3%' 'Active'
$ws.Range("E14").NumberFormat = "@"
Set-RowValues $ws 14 '2026-02-01' '13:49:21' '13:00' 'Bathroom' '82.9%' 'Active'
$ws.Range("E15").NumberFormat = "@"
Set-RowValues $ws 15 '2026-02-01' '13:49:22' '13:00' 'Bathroom' '82.0%' 'Active'
$ws.Range("E16").NumberFormat = "@"
Set-RowValues $ws 16 '2026-02-01' '13:49:23' '13:00' 'Bathroom' '80.7%' 'Active'
$ws.Range("E17").NumberFormat = "@"
Set-RowValues $ws 17 '2026-02-01' '13:49:24' '13:00' 'Bathroom' '81.5%' 'Active'

$ws = $wb.Worksheets.Item("Proximity")
Set-RowValues $ws 20 '2026-02-01' '13:48:40' '13:00' 'Living Room Main Door' 'ENTER' 'User ENTERED Living Room Main Door'
Set-RowValues $ws 21 '2026-02-01' '13:49:18' '13:00' 'Living Room Main Door' 'EXIT' 'User EXITED Living Room Main Door'

$ws = $wb.Worksheets.Item("Camera")
Set-RowValues $ws 13 '2026-02-01' '13:49:17' '13:00' 'Living Room Main Door' 'Image Captured' 'Active'

